# Applies the "add canada bets" update: refreshes the P1-P12 (and Pole)
# predictions for each bettor row and records the newly introduced
# driver code "BEA" (Bearman) in the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("B3").Value = "RUS"
$ws.Range("D3").Value = "NOR"
$ws.Range("E3").Value = "RUS"
$ws.Range("F3").Value = "PIA"
$ws.Range("H3").Value = "VER"
$ws.Range("I3").Value = "ALO"
$ws.Range("J3").Value = "HAM"
$ws.Range("K3").Value = "GAS"
$ws.Range("L3").Value = "OCO"
$ws.Range("M3").Value = "TSU"
$ws.Range("N3").Value = "BEA"
$ws.Range("O3").Value = "HAD"

# Row 4
$ws.Range("B4").Value = "NOR"
$ws.Range("D4").Value = "NOR"
$ws.Range("E4").Value = "VER"
$ws.Range("F4").Value = "RUS"
$ws.Range("G4").Value = "PIA"
$ws.Range("H4").Value = "ANT"
$ws.Range("I4").Value = "HAM"
$ws.Range("J4").Value = "LEC"
$ws.Range("K4").Value = "SAI"
$ws.Range("M4").Value = "ALB"
$ws.Range("N4").Value = "BOR"
$ws.Range("O4").Value = "HAD"

# Row 5
$ws.Range("B5").Value = "VER"
$ws.Range("D5").Value = "VER"
$ws.Range("E5").Value = "SAI"
$ws.Range("F5").Value = "ALB"
$ws.Range("G5").Value = "RUS"
$ws.Range("H5").Value = "HAM"
$ws.Range("I5").Value = "NOR"
$ws.Range("J5").Value = "HAD"
$ws.Range("K5").Value = "GAS"
$ws.Range("L5").Value = "LAW"
$ws.Range("M5").Value = "ALO"
$ws.Range("N5").Value = "PIA"
$ws.Range("O5").Value = "LEC"

# Row 6
$ws.Range("B6").Value = "NOR"
$ws.Range("D6").Value = "NOR"
$ws.Range("E6").Value = "PIA"
$ws.Range("F6").Value = "VER"
$ws.Range("G6").Value = "RUS"
$ws.Range("H6").Value = "ANT"
$ws.Range("I6").Value = "HAM"
$ws.Range("J6").Value = "LEC"
$ws.Range("K6").Value = "ALB"
$ws.Range("L6").Value = "SAI"
$ws.Range("M6").Value = "LAW"
$ws.Range("N6").Value = "HAD"
$ws.Range("O6").Value = "ALO"

# Row 7
$ws.Range("B7").Value = "NOR"
$ws.Range("D7").Value = "NOR"
$ws.Range("E7").Value = "PIA"
$ws.Range("H7").Value = "ANT"
$ws.Range("I7").Value = "LEC"
$ws.Range("K7").Value = "ALB"
$ws.Range("M7").Value = "SAI"
$ws.Range("N7").Value = "HAD"
$ws.Range("O7").Value = "TSU"

# Row 9
$ws.Range("H9").Value = "ANT"
$ws.Range("I9").Value = "LEC"
$ws.Range("J9").Value = "ALB"
$ws.Range("K9").Value = "ALO"
$ws.Range("L9").Value = "GAS"
$ws.Range("M9").Value = "TSU"
$ws.Range("N9").Value = "SAI"
$ws.Range("O9").Value = "HAM"

# Row 10
$ws.Range("B10").Value = "NOR"
$ws.Range("D10").Value = "VER"
$ws.Range("F10").Value = "PIA"
$ws.Range("L10").Value = "SAI"
$ws.Range("M10").Value = "HAD"
$ws.Range("N10").Value = "GAS"
$ws.Range("O10").Value = "ALB"

# Row 11
$ws.Range("B11").Value = "NOR"
$ws.Range("D11").Value = "NOR"
$ws.Range("E11").Value = "VER"
$ws.Range("F11").Value = "PIA"
$ws.Range("K11").Value = "ALB"
$ws.Range("L11").Value = "SAI"
$ws.Range("M11").Value = "HAD"

# Row 12
$ws.Range("E12").Value = "NOR"
$ws.Range("F12").Value = "RUS"
$ws.Range("G12").Value = "VER"
$ws.Range("H12").Value = "ANT"
$ws.Range("I12").Value = "LEC"
$ws.Range("K12").Value = "ALO"
$ws.Range("L12").Value = "OCO"
$ws.Range("N12").Value = "GAS"
$ws.Range("O12").Value = "STR"

# Row 13
$ws.Range("B13").Value = "NOR"
$ws.Range("D13").Value = "NOR"
$ws.Range("E13").Value = "PIA"
$ws.Range("F13").Value = "RUS"
$ws.Range("G13").Value = "VER"
$ws.Range("H13").Value = "ANT"
$ws.Range("J13").Value = "LEC"
$ws.Range("L13").Value = "SAI"
$ws.Range("M13").Value = "ALO"
$ws.Range("N13").Value = "ALB"

# Row 14
$ws.Range("H14").Value = "ANT"
$ws.Range("I14").Value = "LEC"
$ws.Range("K14").Value = "HAD"
$ws.Range("L14").Value = "SAI"
$ws.Range("M14").Value = "ALO"
$ws.Range("N14").Value = "ALB"
$ws.Range("O14").Value = "TSU"

# Row 15
$ws.Range("H15").Value = "ANT"
$ws.Range("I15").Value = "HAM"
$ws.Range("J15").Value = "LEC"
$ws.Range("K15").Value = "ALB"
$ws.Range("L15").Value = "SAI"
$ws.Range("M15").Value = "ALO"
$ws.Range("N15").Value = "HAD"
$ws.Range("O15").Value = "TSU"

# Row 16
$ws.Range("E16").Value = "RUS"
$ws.Range("G16").Value = "ANT"
$ws.Range("H16").Value = "VER"
$ws.Range("J16").Value = "LEC"
$ws.Range("K16").Value = "ALO"
$ws.Range("L16").Value = "ALB"
$ws.Range("M16").Value = "HAD"
$ws.Range("N16").Value = "SAI"

# Row 17
$ws.Range("K17").Value = "ALB"
$ws.Range("L17").Value = "SAI"
$ws.Range("M17").Value = "TSU"
$ws.Range("N17").Value = "HAD"
$ws.Range("O17").Value = "ALO"

# Row 18
$ws.Range("E18").Value = "RUS"
$ws.Range("F18").Value = "PIA"
$ws.Range("G18").Value = "VER"
$ws.Range("H18").Value = "ANT"
$ws.Range("I18").Value = "LEC"
$ws.Range("J18").Value = "TSU"
$ws.Range("K18").Value = "HAM"
$ws.Range("L18").Value = "HAD"
$ws.Range("N18").Value = "HUL"
$ws.Range("O18").Value = "GAS"

# Row 19
$ws.Range("B19").Value = "RUS"
$ws.Range("E19").Value = "RUS"
$ws.Range("G19").Value = "PIA"
$ws.Range("H19").Value = "ANT"
$ws.Range("I19").Value = "ALB"
$ws.Range("J19").Value = "LEC"
$ws.Range("K19").Value = "HAM"
$ws.Range("L19").Value = "HAD"
$ws.Range("M19").Value = "TSU"
$ws.Range("N19").Value = "SAI"
$ws.Range("O19").Value = "BOR"

# Row 20
$ws.Range("B20").Value = "NOR"
$ws.Range("D20").Value = "NOR"
$ws.Range("E20").Value = "RUS"
$ws.Range("F20").Value = "VER"
$ws.Range("G20").Value = "PIA"
$ws.Range("H20").Value = "LEC"
$ws.Range("I20").Value = "HAM"
$ws.Range("J20").Value = "ANT"
$ws.Range("L20").Value = "ALO"
$ws.Range("M20").Value = "TSU"
$ws.Range("N20").Value = "GAS"

# Row 21
$ws.Range("B21").Value = "NOR"
$ws.Range("D21").Value = "NOR"
$ws.Range("F21").Value = "VER"
$ws.Range("G21").Value = "RUS"
$ws.Range("H21").Value = "HAM"
$ws.Range("I21").Value = "ANT"
$ws.Range("J21").Value = "LEC"
$ws.Range("L21").Value = "SAI"
$ws.Range("M21").Value = "TSU"
$ws.Range("N21").Value = "HAD"
$ws.Range("O21").Value = "GAS"

# Row 22
$ws.Range("B22").Value = "NOR"
$ws.Range("D22").Value = "NOR"
$ws.Range("E22").Value = "VER"
$ws.Range("F22").Value = "PIA"
$ws.Range("H22").Value = "HAM"
$ws.Range("I22").Value = "LEC"
$ws.Range("L22").Value = "LAW"
$ws.Range("M22").Value = "ALB"
$ws.Range("N22").Value = "SAI"

# Row 24
$ws.Range("F24").Value = "RUS"
$ws.Range("G24").Value = "VER"
$ws.Range("H24").Value = "LEC"
$ws.Range("K24").Value = "ALB"
$ws.Range("L24").Value = "SAI"
$ws.Range("M24").Value = "HAD"
$ws.Range("N24").Value = "GAS"

# Row 25
$ws.Range("F25").Value = "RUS"
$ws.Range("G25").Value = "VER"
$ws.Range("M25").Value = "HAD"
$ws.Range("N25").Value = "TSU"
$ws.Range("O25").Value = "ALO"

# Row 26
$ws.Range("E26").Value = "PIA"
$ws.Range("F26").Value = "VER"
$ws.Range("G26").Value = "LEC"
$ws.Range("J26").Value = "HAD"
$ws.Range("K26").Value = "ALB"
$ws.Range("L26").Value = "SAI"
$ws.Range("M26").Value = "HUL"
$ws.Range("N26").Value = "STR"
$ws.Range("O26").Value = "ANT"

# Row 27
$ws.Range("B27").Value = "PIA"
$ws.Range("G27").Value = "RUS"
$ws.Range("H27").Value = "HAM"
$ws.Range("J27").Value = "ANT"
$ws.Range("K27").Value = "LAW"
$ws.Range("L27").Value = "ALO"
$ws.Range("M27").Value = "SAI"
$ws.Range("N27").Value = "GAS"
$ws.Range("O27").Value = "TSU"

# Restore the cursor/selection position recorded in the workbook
$null = $ws.Range("Q14").Select()
